$wb = $excel.ActiveWorkbook

# --- Sheet "Input": D3 Vrouw -> Man, selection -> F5 ---
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("D3").Value = "Man"
$wsInput.Activate()
$wsInput.Range("F5").Select() | Out-Null

# --- Sheet "Uitgaven": new Woonlasten row 3, selection -> D4 ---
$wsUitgaven = $wb.Worksheets.Item("Uitgaven")
$wsUitgaven.Range("A3").Value = "Woonlasten"
$wsUitgaven.Range("B3").Value = 450
$wsUitgaven.Range("C3").Value = "15-06-2012"
$wsUitgaven.Range("D3").Value = "30-09-2053"
$wsUitgaven.Activate()
$wsUitgaven.Range("D4").Select() | Out-Null

# --- Sheet "Output": B4 formula changed, selection -> B5 ---
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("B4").Formula = "=YEAR('Input Oud'!B19)-YEAR('Input Oud'!B4)+(MONTH('Input Oud'!B19)-MONTH('Input Oud'!B4))/12"
$wsOutput.Activate()
$wsOutput.Range("B5").Select() | Out-Null

# --- Sheet "Input Oud": B6 formula changed ---
$wsInputOud = $wb.Worksheets.Item("Input Oud")
$wsInputOud.Range("B6").Formula = "=2016-YEAR(B4)+(8-MONTH(B4))/12+(1-DAY(B4))/365"

# --- Restore active sheet to "Input" (tabSelected) ---
$wsInput.Activate()
